$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# New export rows for 2025-11-05 and 2025-11-06.
# Format the Date column as text first so the ISO-looking strings are
# stored as plain text (matching the rest of the Date column) instead of
# being auto-converted to date serials, then drop back to the default
# (unstyled) look once the values are in place.
$ws.Range("A32:A33").NumberFormat = "@"

$ws.Range("A32").Value = "2025-11-05"
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 101

$ws.Range("A33").Value = "2025-11-06"
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 95

$ws.Range("A32:A33").ClearFormats()
